# Auto-generated Excel COM-interop script implementing the row-content re-sort
# for rows 21-27 and 55-56 on the 'Artfynd' worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = 131047027
$ws.Range("B21").Value = 78909
$ws.Range("E21").Value = 353
$ws.Range("F21").Value = 'Dvärgbägarlav'
$ws.Range("G21").Value = 'Cladonia parasitica'
$ws.Range("H21").Value = '(Hoffm.) Hoffm.'
$ws.Range("Q21").Value = 395391
$ws.Range("R21").Value = 6804697
$ws.Range("Z21").Value = '11:50'
$ws.Range("AB21").Value = '11:50'
$ws.Range("A22").Value = 131046963
$ws.Range("B22").Value = 79243
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = 'Garnlav'
$ws.Range("G22").Value = 'Alectoria sarmentosa'
$ws.Range("H22").Value = '(Ach.) Ach.'
$ws.Range("Q22").Value = 395386
$ws.Range("R22").Value = 6804723
$ws.Range("Z22").Value = '11:52'
$ws.Range("AB22").Value = '11:52'
$ws.Range("A23").Value = 131046972
$ws.Range("Q23").Value = 395524
$ws.Range("R23").Value = 6804798
$ws.Range("Z23").Value = '12:03'
$ws.Range("AB23").Value = '12:03'
$ws.Range("A24").Value = 131046792
$ws.Range("B24").Value = 57884
$ws.Range("E24").Value = 100109
$ws.Range("F24").Value = 'Tretåig hackspett'
$ws.Range("G24").Value = 'Picoides tridactylus'
$ws.Range("H24").Value = '(Linnaeus, 1758)'
$ws.Range("M24").Value = 'färska spår'
$ws.Range("Q24").Value = 395357
$ws.Range("R24").Value = 6804769
$ws.Range("Z24").Value = '11:22'
$ws.Range("AB24").Value = '11:22'
$ws.Range("AC24").Value = 'Färska ringhack (gran)'
$ws.Range("A25").Value = 131047020
$ws.Range("B25").Value = 57884
$ws.Range("E25").Value = 100109
$ws.Range("F25").Value = 'Tretåig hackspett'
$ws.Range("G25").Value = 'Picoides tridactylus'
$ws.Range("H25").Value = '(Linnaeus, 1758)'
$ws.Range("M25").Value = 'färska spår'
$ws.Range("Q25").Value = 395541
$ws.Range("R25").Value = 6804800
$ws.Range("Z25").Value = '12:04'
$ws.Range("AB25").Value = '12:04'
$ws.Range("AC25").Value = 'Troliga spår efter tretåig hackspett (barkfälkning)'
$ws.Range("AE25").Value = $true
$ws.Range("A26").Value = 131046714
$ws.Range("B26").Value = 83223
$ws.Range("E26").Value = 6440
$ws.Range("F26").Value = 'Vitgrynig nållav'
$ws.Range("G26").Value = 'Chaenotheca subroscida'
$ws.Range("H26").Value = '(Eitner) Zahlbr.'
$ws.Range("M26").Value = ""
$ws.Range("Q26").Value = 395419
$ws.Range("R26").Value = 6804801
$ws.Range("Z26").Value = '11:27'
$ws.Range("AB26").Value = '11:27'
$ws.Range("AC26").Value = ""
$ws.Range("A27").Value = 131046922
$ws.Range("B27").Value = 79243
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = 'Garnlav'
$ws.Range("G27").Value = 'Alectoria sarmentosa'
$ws.Range("H27").Value = '(Ach.) Ach.'
$ws.Range("M27").Value = ""
$ws.Range("Q27").Value = 395365
$ws.Range("R27").Value = 6804755
$ws.Range("Z27").Value = '11:20'
$ws.Range("AB27").Value = '11:20'
$ws.Range("AC27").Value = ""
$ws.Range("AE27").Value = $false
$ws.Range("A55").Value = 131046925
$ws.Range("B55").Value = 79243
$ws.Range("E55").Value = 6425
$ws.Range("F55").Value = 'Garnlav'
$ws.Range("G55").Value = 'Alectoria sarmentosa'
$ws.Range("H55").Value = '(Ach.) Ach.'
$ws.Range("Q55").Value = 395380
$ws.Range("R55").Value = 6804774
$ws.Range("Z55").Value = '11:25'
$ws.Range("AB55").Value = '11:25'
$ws.Range("A56").Value = 131046722
$ws.Range("B56").Value = 79275
$ws.Range("E56").Value = 185
$ws.Range("F56").Value = 'Violettgrå tagellav'
$ws.Range("G56").Value = 'Bryoria nadvornikiana'
$ws.Range("H56").Value = '(Gyeln.) Brodo & D.Hawksw.'
$ws.Range("Q56").Value = 395391
$ws.Range("R56").Value = 6804603
$ws.Range("Z56").Value = '10:52'
$ws.Range("AB56").Value = '10:52'
